# Minor re-balancing of weapons: make strike craft even harder to hit
# by lowering the "Vs. strike craft" (column N) values on the Warheads sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Warheads")

# Map of row -> new value for column N ("Vs. strike craft")
$updates = @{
    2  = 0.3
    3  = 0.3
    4  = 0.3
    6  = 0.2
    7  = 0.2
    8  = 0.5
    10 = 0.2
    11 = 0.2
    12 = 0.5
    20 = 0.2
    21 = 0.2
    22 = 0.4
    24 = 0.15
    25 = 0.15
    26 = 0.5
    28 = 0.15
    29 = 0.15
    30 = 0.5
    38 = 0.15
    39 = 0.15
    40 = 0.65
    44 = 0.75
    48 = 0.75
}

foreach ($row in $updates.Keys) {
    $ws.Range("N$row").Value = $updates[$row]
}

# Restore the view so the top-left cell is back at A1 and the selection
# is set to E5 (matching the saved view state in the workbook).
$ws.Activate()
$ws.Range("A1").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E5").Select() | Out-Null
